$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 21 new rows at row 15 (pushes the existing "Chorus 1..81" block,
# currently at rows 15-95, down to rows 36-116) to make room for a new
# "Bridge 1..21" section.
$ws.Range("A15:A35").EntireRow.Insert()

$bridgeData = @(
    @("Bridge 1", "Like I'm the only one that's in command"),
    @("Bridge 2", "Only one"),
    @("Bridge 3", "Only one"),
    @("Bridge 4", "Take me for a ride, ride"),
    @("Bridge 5", "Oh, baby, take me high, high"),
    @("Bridge 6", "Let me make you rise, rise"),
    @("Bridge 7", "Take me for a ride, ride"),
    @("Bridge 8", "Oh, baby, take me high, high"),
    @("Bridge 9", "Let me make you rise, rise"),
    @("Bridge 10", "Oh, baby, take me high, high"),
    @("Bridge 11", "Let me make you rise, rise"),
    @("Bridge 12", "Let me make you rise, rise"),
    @("Bridge 13", "Take me for a ride, ride"),
    @("Bridge 14", "Oh, baby, take me high, high"),
    @("Bridge 15", "Let me make you rise, rise"),
    @("Bridge 16", "Oh, baby, take me high, high"),
    @("Bridge 17", "Let me make you rise, rise"),
    @("Bridge 18", "Let me make you rise, rise"),
    @("Bridge 19", "Like I'm the only one that's in command"),
    @("Bridge 20", "Girl in the world"),
    @("Bridge 21", "Girl in the world")
)

$row = 15
foreach ($entry in $bridgeData) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
